$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> updated values for columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$rows = @{
    2  = @{ D = 44874; L = "Especial"; M = 30;  N = 25000; O = 25000; P = 25000; S = 2500 }
    3  = @{ D = 44874; L = "Primera";  M = 80;  N = 23000; O = 24000; P = 23500; S = 2350 }
    4  = @{ D = 44447; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 }
    5  = @{ D = 44461; L = "Especial"; M = 60;  N = 31000; O = 32000; P = 31500; S = 3150 }
    6  = @{ D = 44461; L = "Primera";  M = 30;  N = 30000; O = 30000; P = 30000; S = 3000 }
    7  = @{ D = 44446; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 }
    8  = @{ D = 44460; L = "Especial"; M = 60;  N = 31000; O = 32000; P = 31500; S = 3150 }
    9  = @{ D = 44460; L = "Primera";  M = 30;  N = 30000; O = 30000; P = 30000; S = 3000 }
    10 = @{ D = 44839; L = "Primera";  M = 120; N = 25000; O = 26000; P = 25500; S = 2550 }
    11 = @{ D = 44448; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 }
    12 = @{ D = 44868; L = "Especial"; M = 60;  N = 26000; O = 26000; P = 26000; S = 2600 }
    13 = @{ D = 44452; L = "Primera";  M = 60;  N = 21000; O = 22000; P = 21500; S = 2150 }
    14 = @{ D = 44848; L = "Especial"; M = 60;  N = 24000; O = 25000; P = 24500; S = 2450 }
    15 = @{ D = 44848; L = "Primera";  M = 120; N = 21000; O = 22000; P = 21500; S = 2150 }
    16 = @{ D = 44487; L = "Primera";  M = 30;  N = 23000; O = 24000; P = 23500; S = 2350 }
    17 = @{ D = 44841; L = "Primera";  M = 60;  N = 23000; O = 24000; P = 23500; S = 2350 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value2 = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("S$r").Value = $vals.S
}
